$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price records in rows 2, 3 and 5 are rotated: the data that
# used to live in row 5 now belongs to row 2, the data that used to live
# in row 2 now belongs to row 3, and the data that used to live in row 3
# now belongs to row 5 (row 4 is unaffected).

# --- Row 2 gets what used to be in row 5 ---
$ws.Range("D2").Value = 44623
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1900
$ws.Range("N2").Value = "`$/paquete"
$ws.Range("P2").Value = 1900
$ws.Range("Q2").Value = 1

# --- Row 3 gets what used to be in row 2 ---
$ws.Range("D3").Value = 44370
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1200
$ws.Range("M3").Value = 1080
$ws.Range("N3").Value = "`$/docena de matas"
$ws.Range("P3").Value = 180
$ws.Range("Q3").Value = 6

# --- Row 5 gets what used to be in row 3 ---
$ws.Range("D5").Value = 44377
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2800
$ws.Range("M5").Value = 2364
$ws.Range("N5").Value = "`$/docena de matas"
$ws.Range("P5").Value = 394
$ws.Range("Q5").Value = 6
